$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 175, shifting existing rows 175:268 down to 176:269.
$ws.Rows(175).Insert()

# Populate the newly inserted row 175 with its data (matches the row that
# used to be at 175 for most columns, with updated D/M/N/O/P/S values).
$ws.Range("A175").Value = 10
$ws.Range("B175").Value = "Vega Modelo de Temuco"
$ws.Range("C175").Value = "La Araucanía"
$ws.Range("D175").Value = 44452
$ws.Range("E175").Value = 9
$ws.Range("F175").Value = "Fruta"
$ws.Range("G175").Value = 100108
$ws.Range("H175").Value = "Tropicales y subtropicales"
$ws.Range("I175").Value = 100108005
$ws.Range("J175").Value = "Piña"
$ws.Range("K175").Value = "Caramelo"
$ws.Range("L175").Value = "Primera"
$ws.Range("M175").Value = 110
$ws.Range("N175").Value = 22000
$ws.Range("O175").Value = 22000
$ws.Range("P175").Value = 22000
$ws.Range("Q175").Value = "$/caja 12 unidades"
$ws.Range("R175").Value = "Ecuador"
$ws.Range("S175").Value = 1833
$ws.Range("T175").Value = 12
